# feat: add 2022-Q3 data
#
# The workbook previously had two sheets: "总计" (totals) and "2022-Q2"
# (fund detail for that quarter). This script adds a new "2022-Q3" sheet
# with that quarter's fund detail, keeps the original "2022-Q2" detail
# sheet intact (moved to the end), and updates the "总计" summary sheet
# with a new row for 2022-Q3 while pushing the old 2022-Q2 summary row down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet so its data survives
# under its own tab again; the original sheet object is renamed/reused
# below to become the new "2022-Q3" sheet.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)

$dup = $wb.Worksheets.Item("2022-Q2 (2)")
$dup.Name = "2022-Q2-staging"
$q2.Name = "2022-Q3"
$dup.Name = "2022-Q2"

# ---------------------------------------------------------------------
# Step 2: populate "2022-Q3" with this quarter's fund data (only two
# funds this time, so the old third data row is dropped).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Range("A4:H4").Clear()

# Fund 1: 519029 / 华夏稳增混合
$q3.Range("B2").Value = "'519029"
$q3.Range("B2").Style = "Normal"
$q3.Range("C2").Value = "华夏稳增混合"
$q3.Range("D2").Value = "'8.56"
$q3.Range("D2").Style = "Normal"
$q3.Range("E2").Value = "'93.73"
$q3.Range("E2").Style = "Normal"
$q3.Range("F2").Value = "'5.80"
$q3.Range("F2").Style = "Normal"
$q3.Range("G2").Value = "'0.4965"
$q3.Range("G2").Style = "Normal"
$q3.Range("H2").Value = 3

# Fund 2: 161040 / 富国创业板两年定期开放混合
$q3.Range("B3").Value = "'161040"
$q3.Range("B3").Style = "Normal"
$q3.Range("C3").Value = "富国创业板两年定期开放混合"
$q3.Range("D3").Value = "'11.31"
$q3.Range("D3").Style = "Normal"
$q3.Range("E3").Value = "'70.93"
$q3.Range("E3").Style = "Normal"
$q3.Range("F3").Value = "'3.33"
$q3.Range("F3").Style = "Normal"
$q3.Range("G3").Value = "'0.3766"
$q3.Range("G3").Style = "Normal"
$q3.Range("H3").Value = 7

# Re-apply this workbook's header / first-column style (the bold,
# bordered, centred look used by the "总计" sheet) to the new sheet.
$total = $wb.Worksheets.Item("总计")
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - push the old 2022-Q2 row
# down to row 3 and write the new 2022-Q3 row into row 2.
# ---------------------------------------------------------------------
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.67
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.87

# Restore "总计" as the active sheet/tab.
$total.Activate()
